$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 27780988
$ws.Cells.Item(100, 9).Value = 41669230
$ws.Cells.Item(100, 10).Value = 4500
$ws.Cells.Item(100, 11).Value = 41669230
$ws.Cells.Item(100, 12).Value = 4500
$ws.Cells.Item(100, 13).Value = -41668689
$ws.Cells.Item(100, 14).Value = -5582
$ws.Cells.Item(138, 8).Value = 4377.72
$ws.Cells.Item(138, 9).Value = 1207.591
$ws.Cells.Item(138, 10).Value = 6868.5356
$ws.Cells.Item(138, 11).Value = 3622.773
$ws.Cells.Item(138, 12).Value = 20605.6068
$ws.Cells.Item(138, 13).Value = 1517.227
$ws.Cells.Item(138, 14).Value = -30885.6068

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(27, 8).Value = 3695.8
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 3695.8
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 3695.8
$ws.Cells.Item(27, 13).Value = $null
$ws.Cells.Item(27, 14).Value = -4063.8
$ws.Cells.Item(40, 8).Value = 4600
$ws.Cells.Item(40, 10).Value = 4600
$ws.Cells.Item(40, 12).Value = 4600
$ws.Cells.Item(40, 14).Value = -4952
$ws.Cells.Item(58, 8).Value = 5353
$ws.Cells.Item(58, 9).Value = 5353
$ws.Cells.Item(58, 11).Value = 5353
$ws.Cells.Item(58, 13).Value = -4923
$ws.Cells.Item(110, 8).Value = 979.6429000000001
$ws.Cells.Item(110, 9).Value = 840.25
$ws.Cells.Item(110, 10).Value = 1165.5
$ws.Cells.Item(110, 11).Value = 840.25
$ws.Cells.Item(110, 12).Value = 1165.5
$ws.Cells.Item(110, 13).Value = 1204.75
$ws.Cells.Item(110, 14).Value = -5255.5
$ws.Cells.Item(132, 8).Value = 19251512
$ws.Cells.Item(132, 9).Value = 23256830
$ws.Cells.Item(132, 10).Value = 114990.336
$ws.Cells.Item(132, 11).Value = 69770490
$ws.Cells.Item(132, 12).Value = 344971.008
$ws.Cells.Item(132, 13).Value = -69767960
$ws.Cells.Item(132, 14).Value = -350031.008

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(57, 8).Value = 8461
$ws.Cells.Item(57, 10).Value = 8461
$ws.Cells.Item(57, 12).Value = 8461
$ws.Cells.Item(57, 14).Value = -9581
$ws.Cells.Item(80, 8).Value = 35000
$ws.Cells.Item(80, 10).Value = 35000
$ws.Cells.Item(80, 12).Value = 35000
$ws.Cells.Item(80, 14).Value = -37246
$ws.Cells.Item(83, 8).Value = 35000
$ws.Cells.Item(83, 10).Value = 35000
$ws.Cells.Item(83, 12).Value = 105000
$ws.Cells.Item(83, 14).Value = -116232
$ws.Cells.Item(132, 8).Value = 26621.75
$ws.Cells.Item(132, 9).Value = 1324.8857
$ws.Cells.Item(132, 10).Value = 203699.8
$ws.Cells.Item(132, 11).Value = 3974.6571
$ws.Cells.Item(132, 12).Value = 611099.3999999999
$ws.Cells.Item(132, 13).Value = -1444.6571
$ws.Cells.Item(132, 14).Value = -616159.3999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 891.8333
$ws.Cells.Item(5, 9).Value = 826.64703
$ws.Cells.Item(5, 11).Value = 2479.94109
$ws.Cells.Item(5, 13).Value = -2367.94109
$ws.Cells.Item(113, 8).Value = 946.04346
$ws.Cells.Item(113, 9).Value = 684.5
$ws.Cells.Item(113, 10).Value = 985.275
$ws.Cells.Item(113, 11).Value = 2053.5
$ws.Cells.Item(113, 12).Value = 2955.825
$ws.Cells.Item(113, 13).Value = 116.5
$ws.Cells.Item(113, 14).Value = -7295.825
$ws.Cells.Item(135, 8).Value = 891.8333
$ws.Cells.Item(135, 9).Value = 826.64703
$ws.Cells.Item(135, 11).Value = 7439.82327
$ws.Cells.Item(135, 13).Value = -4904.82327

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 285716400
$ws.Cells.Item(33, 10).Value = 285716400
$ws.Cells.Item(33, 12).Value = 285716400
$ws.Cells.Item(33, 14).Value = -285716904
$ws.Cells.Item(36, 8).Value = 2333.3333
$ws.Cells.Item(36, 9).Value = 1000
$ws.Cells.Item(36, 10).Value = 5000
$ws.Cells.Item(36, 11).Value = 1000
$ws.Cells.Item(36, 12).Value = 5000
$ws.Cells.Item(36, 13).Value = -515
$ws.Cells.Item(36, 14).Value = -5970
$ws.Cells.Item(40, 8).Value = 3575
$ws.Cells.Item(40, 10).Value = 3575
$ws.Cells.Item(40, 12).Value = 3575
$ws.Cells.Item(40, 14).Value = -3877
$ws.Cells.Item(43, 8).Value = 16320
$ws.Cells.Item(43, 10).Value = 19900
$ws.Cells.Item(43, 12).Value = 19900
$ws.Cells.Item(43, 14).Value = -20202
$ws.Cells.Item(132, 8).Value = 38749.184
$ws.Cells.Item(132, 9).Value = 1418.25
$ws.Cells.Item(132, 11).Value = 4254.75
$ws.Cells.Item(132, 13).Value = -1724.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(34, 8).Value = 3998.5
$ws.Cells.Item(34, 9).Value = 4000
$ws.Cells.Item(34, 10).Value = 3998.2
$ws.Cells.Item(34, 11).Value = 4000
$ws.Cells.Item(34, 12).Value = 3998.2
$ws.Cells.Item(34, 13).Value = -3828
$ws.Cells.Item(34, 14).Value = -4342.2
$ws.Cells.Item(39, 8).Value = 3825
$ws.Cells.Item(39, 9).Value = 2500
$ws.Cells.Item(39, 10).Value = 4266.6665
$ws.Cells.Item(39, 11).Value = 2500
$ws.Cells.Item(39, 12).Value = 4266.6665
$ws.Cells.Item(39, 13).Value = -2040
$ws.Cells.Item(39, 14).Value = -5186.6665
$ws.Cells.Item(45, 8).Value = 5900
$ws.Cells.Item(45, 10).Value = 5900
$ws.Cells.Item(45, 12).Value = 5900
$ws.Cells.Item(45, 14).Value = -6714
$ws.Cells.Item(47, 8).Value = 5231.5
$ws.Cells.Item(47, 9).Value = 4059
$ws.Cells.Item(47, 11).Value = 4059
$ws.Cells.Item(47, 13).Value = -3569
$ws.Cells.Item(52, 8).Value = 5231.5
$ws.Cells.Item(52, 9).Value = 4059
$ws.Cells.Item(52, 11).Value = 4059
$ws.Cells.Item(52, 13).Value = -3826
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).Value = $null
$ws.Cells.Item(56, 8).Value = 15650.333
$ws.Cells.Item(56, 9).Value = 10051
$ws.Cells.Item(56, 10).Value = 18450
$ws.Cells.Item(56, 11).Value = 10051
$ws.Cells.Item(56, 12).Value = 18450
$ws.Cells.Item(56, 13).Value = -9360
$ws.Cells.Item(56, 14).Value = -19832
$ws.Cells.Item(70, 8).Value = 7750
$ws.Cells.Item(70, 9).Value = 10000
$ws.Cells.Item(70, 10).Value = 5500
$ws.Cells.Item(70, 11).Value = 10000
$ws.Cells.Item(70, 12).Value = 5500
$ws.Cells.Item(70, 13).Value = -9730
$ws.Cells.Item(70, 14).Value = -6040
$ws.Cells.Item(73, 8).Value = 7750
$ws.Cells.Item(73, 9).Value = 10000
$ws.Cells.Item(73, 10).Value = 5500
$ws.Cells.Item(73, 11).Value = 10000
$ws.Cells.Item(73, 12).Value = 5500
$ws.Cells.Item(73, 13).Value = -9064
$ws.Cells.Item(73, 14).Value = -7372
$ws.Cells.Item(74, 8).Value = 13291
$ws.Cells.Item(74, 9).Value = 11618.333
$ws.Cells.Item(74, 10).Value = 15800
$ws.Cells.Item(74, 11).Value = 11618.333
$ws.Cells.Item(74, 12).Value = 15800
$ws.Cells.Item(74, 13).Value = -10620.333
$ws.Cells.Item(74, 14).Value = -17796
$ws.Cells.Item(77, 8).Value = 13291
$ws.Cells.Item(77, 9).Value = 11618.333
$ws.Cells.Item(77, 10).Value = 15800
$ws.Cells.Item(77, 11).Value = 34854.999
$ws.Cells.Item(77, 12).Value = 47400
$ws.Cells.Item(77, 13).Value = -29862.999
$ws.Cells.Item(133, 8).Value = 102663.336
$ws.Cells.Item(133, 10).Value = 102663.336
$ws.Cells.Item(133, 12).Value = 102663.336
$ws.Cells.Item(133, 14).Value = -107723.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(50, 8).Value = 6400
$ws.Cells.Item(50, 10).Value = 6400
$ws.Cells.Item(50, 12).Value = 6400
$ws.Cells.Item(50, 14).Value = -7662
$ws.Cells.Item(58, 8).Value = 3925.5715
$ws.Cells.Item(58, 9).Value = 2369.75
$ws.Cells.Item(58, 11).Value = 2369.75
$ws.Cells.Item(58, 13).Value = -2061.75
$ws.Cells.Item(64, 8).Value = 6000
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 13).Value = $null
$ws.Cells.Item(67, 8).Value = 6000
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 13).Value = $null
$ws.Cells.Item(82, 8).Value = 32000
$ws.Cells.Item(82, 10).Value = 32000
$ws.Cells.Item(82, 12).Value = 32000
$ws.Cells.Item(82, 14).Value = -32766
$ws.Cells.Item(85, 8).Value = 32000
$ws.Cells.Item(85, 10).Value = 32000
$ws.Cells.Item(85, 12).Value = 32000
$ws.Cells.Item(85, 14).Value = -34652
$ws.Cells.Item(107, 8).Value = 708.3889
$ws.Cells.Item(107, 9).Value = 717.3333
$ws.Cells.Item(107, 10).Value = 703.9167
$ws.Cells.Item(107, 11).Value = 2151.9999
$ws.Cells.Item(107, 12).Value = 2111.7501
$ws.Cells.Item(107, 13).Value = -231.9998999999998
$ws.Cells.Item(107, 14).Value = -5951.7501
$ws.Cells.Item(132, 8).Value = 85888424
$ws.Cells.Item(132, 9).Value = 150001390
$ws.Cells.Item(132, 11).Value = 450004170
$ws.Cells.Item(132, 13).Value = -450001640
